$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the "species record" block (columns A,B,D,E,F,G,H,M,Q,R,AO) for each
# source row BEFORE writing anything, since the edit is a single 10-cycle
# permutation across rows 4,5,6,7,8,10,11,12,13,14 (row 9 is untouched).
$A4 = $ws.Range("A4").Value2
$B4 = $ws.Range("B4").Value2
$D4 = $ws.Range("D4").Value2
$E4 = $ws.Range("E4").Value2
$F4 = $ws.Range("F4").Value2
$G4 = $ws.Range("G4").Value2
$H4 = $ws.Range("H4").Value2
$M4 = $ws.Range("M4").Value2
$Q4 = $ws.Range("Q4").Value2
$R4 = $ws.Range("R4").Value2
$AO4 = $ws.Range("AO4").Value2
$A5 = $ws.Range("A5").Value2
$B5 = $ws.Range("B5").Value2
$D5 = $ws.Range("D5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$M5 = $ws.Range("M5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2
$AO5 = $ws.Range("AO5").Value2
$A6 = $ws.Range("A6").Value2
$B6 = $ws.Range("B6").Value2
$D6 = $ws.Range("D6").Value2
$E6 = $ws.Range("E6").Value2
$F6 = $ws.Range("F6").Value2
$G6 = $ws.Range("G6").Value2
$H6 = $ws.Range("H6").Value2
$M6 = $ws.Range("M6").Value2
$Q6 = $ws.Range("Q6").Value2
$R6 = $ws.Range("R6").Value2
$AO6 = $ws.Range("AO6").Value2
$A7 = $ws.Range("A7").Value2
$B7 = $ws.Range("B7").Value2
$D7 = $ws.Range("D7").Value2
$E7 = $ws.Range("E7").Value2
$F7 = $ws.Range("F7").Value2
$G7 = $ws.Range("G7").Value2
$H7 = $ws.Range("H7").Value2
$M7 = $ws.Range("M7").Value2
$Q7 = $ws.Range("Q7").Value2
$R7 = $ws.Range("R7").Value2
$AO7 = $ws.Range("AO7").Value2
$A8 = $ws.Range("A8").Value2
$B8 = $ws.Range("B8").Value2
$D8 = $ws.Range("D8").Value2
$E8 = $ws.Range("E8").Value2
$F8 = $ws.Range("F8").Value2
$G8 = $ws.Range("G8").Value2
$H8 = $ws.Range("H8").Value2
$M8 = $ws.Range("M8").Value2
$Q8 = $ws.Range("Q8").Value2
$R8 = $ws.Range("R8").Value2
$AO8 = $ws.Range("AO8").Value2
$A10 = $ws.Range("A10").Value2
$B10 = $ws.Range("B10").Value2
$D10 = $ws.Range("D10").Value2
$E10 = $ws.Range("E10").Value2
$F10 = $ws.Range("F10").Value2
$G10 = $ws.Range("G10").Value2
$H10 = $ws.Range("H10").Value2
$M10 = $ws.Range("M10").Value2
$Q10 = $ws.Range("Q10").Value2
$R10 = $ws.Range("R10").Value2
$AO10 = $ws.Range("AO10").Value2
$A11 = $ws.Range("A11").Value2
$B11 = $ws.Range("B11").Value2
$D11 = $ws.Range("D11").Value2
$E11 = $ws.Range("E11").Value2
$F11 = $ws.Range("F11").Value2
$G11 = $ws.Range("G11").Value2
$H11 = $ws.Range("H11").Value2
$M11 = $ws.Range("M11").Value2
$Q11 = $ws.Range("Q11").Value2
$R11 = $ws.Range("R11").Value2
$AO11 = $ws.Range("AO11").Value2
$A12 = $ws.Range("A12").Value2
$B12 = $ws.Range("B12").Value2
$D12 = $ws.Range("D12").Value2
$E12 = $ws.Range("E12").Value2
$F12 = $ws.Range("F12").Value2
$G12 = $ws.Range("G12").Value2
$H12 = $ws.Range("H12").Value2
$M12 = $ws.Range("M12").Value2
$Q12 = $ws.Range("Q12").Value2
$R12 = $ws.Range("R12").Value2
$AO12 = $ws.Range("AO12").Value2
$A13 = $ws.Range("A13").Value2
$B13 = $ws.Range("B13").Value2
$D13 = $ws.Range("D13").Value2
$E13 = $ws.Range("E13").Value2
$F13 = $ws.Range("F13").Value2
$G13 = $ws.Range("G13").Value2
$H13 = $ws.Range("H13").Value2
$M13 = $ws.Range("M13").Value2
$Q13 = $ws.Range("Q13").Value2
$R13 = $ws.Range("R13").Value2
$AO13 = $ws.Range("AO13").Value2
$A14 = $ws.Range("A14").Value2
$B14 = $ws.Range("B14").Value2
$D14 = $ws.Range("D14").Value2
$E14 = $ws.Range("E14").Value2
$F14 = $ws.Range("F14").Value2
$G14 = $ws.Range("G14").Value2
$H14 = $ws.Range("H14").Value2
$M14 = $ws.Range("M14").Value2
$Q14 = $ws.Range("Q14").Value2
$R14 = $ws.Range("R14").Value2
$AO14 = $ws.Range("AO14").Value2

# Now write each destination row with the values captured from its source row.
# row 4 <- old row 8
$ws.Range("A4").Value = $A8
$ws.Range("B4").Value = $B8
$ws.Range("D4").Value = $D8
$ws.Range("E4").Value = $E8
$ws.Range("F4").Value = $F8
$ws.Range("G4").Value = $G8
$ws.Range("H4").Value = $H8
$ws.Range("M4").Value = $M8
$ws.Range("Q4").Value = $Q8
$ws.Range("R4").Value = $R8
$ws.Range("AO4").Value = $AO8

# row 5 <- old row 13
$ws.Range("A5").Value = $A13
$ws.Range("B5").Value = $B13
$ws.Range("D5").Value = $D13
$ws.Range("E5").Value = $E13
$ws.Range("F5").Value = $F13
$ws.Range("G5").Value = $G13
$ws.Range("H5").Value = $H13
$ws.Range("M5").Value = $M13
$ws.Range("Q5").Value = $Q13
$ws.Range("R5").Value = $R13
$ws.Range("AO5").Value = $AO13

# row 6 <- old row 5
$ws.Range("A6").Value = $A5
$ws.Range("B6").Value = $B5
$ws.Range("D6").Value = $D5
$ws.Range("E6").Value = $E5
$ws.Range("F6").Value = $F5
$ws.Range("G6").Value = $G5
$ws.Range("H6").Value = $H5
$ws.Range("M6").Value = $M5
$ws.Range("Q6").Value = $Q5
$ws.Range("R6").Value = $R5
$ws.Range("AO6").Value = $AO5

# row 7 <- old row 4
$ws.Range("A7").Value = $A4
$ws.Range("B7").Value = $B4
$ws.Range("D7").Value = $D4
$ws.Range("E7").Value = $E4
$ws.Range("F7").Value = $F4
$ws.Range("G7").Value = $G4
$ws.Range("H7").Value = $H4
$ws.Range("M7").Value = $M4
$ws.Range("Q7").Value = $Q4
$ws.Range("R7").Value = $R4
$ws.Range("AO7").Value = $AO4

# row 8 <- old row 11
$ws.Range("A8").Value = $A11
$ws.Range("B8").Value = $B11
$ws.Range("D8").Value = $D11
$ws.Range("E8").Value = $E11
$ws.Range("F8").Value = $F11
$ws.Range("G8").Value = $G11
$ws.Range("H8").Value = $H11
$ws.Range("M8").Value = $M11
$ws.Range("Q8").Value = $Q11
$ws.Range("R8").Value = $R11
$ws.Range("AO8").Value = $AO11

# row 10 <- old row 14
$ws.Range("A10").Value = $A14
$ws.Range("B10").Value = $B14
$ws.Range("D10").Value = $D14
$ws.Range("E10").Value = $E14
$ws.Range("F10").Value = $F14
$ws.Range("G10").Value = $G14
$ws.Range("H10").Value = $H14
$ws.Range("M10").Value = $M14
$ws.Range("Q10").Value = $Q14
$ws.Range("R10").Value = $R14
$ws.Range("AO10").Value = $AO14

# row 11 <- old row 10
$ws.Range("A11").Value = $A10
$ws.Range("B11").Value = $B10
$ws.Range("D11").Value = $D10
$ws.Range("E11").Value = $E10
$ws.Range("F11").Value = $F10
$ws.Range("G11").Value = $G10
$ws.Range("H11").Value = $H10
$ws.Range("M11").Value = $M10
$ws.Range("Q11").Value = $Q10
$ws.Range("R11").Value = $R10
$ws.Range("AO11").Value = $AO10

# row 12 <- old row 6
$ws.Range("A12").Value = $A6
$ws.Range("B12").Value = $B6
$ws.Range("D12").Value = $D6
$ws.Range("E12").Value = $E6
$ws.Range("F12").Value = $F6
$ws.Range("G12").Value = $G6
$ws.Range("H12").Value = $H6
$ws.Range("M12").Value = $M6
$ws.Range("Q12").Value = $Q6
$ws.Range("R12").Value = $R6
$ws.Range("AO12").Value = $AO6

# row 13 <- old row 7
$ws.Range("A13").Value = $A7
$ws.Range("B13").Value = $B7
$ws.Range("D13").Value = $D7
$ws.Range("E13").Value = $E7
$ws.Range("F13").Value = $F7
$ws.Range("G13").Value = $G7
$ws.Range("H13").Value = $H7
$ws.Range("M13").Value = $M7
$ws.Range("Q13").Value = $Q7
$ws.Range("R13").Value = $R7
$ws.Range("AO13").Value = $AO7

# row 14 <- old row 12
$ws.Range("A14").Value = $A12
$ws.Range("B14").Value = $B12
$ws.Range("D14").Value = $D12
$ws.Range("E14").Value = $E12
$ws.Range("F14").Value = $F12
$ws.Range("G14").Value = $G12
$ws.Range("H14").Value = $H12
$ws.Range("M14").Value = $M12
$ws.Range("Q14").Value = $Q12
$ws.Range("R14").Value = $R12
$ws.Range("AO14").Value = $AO12
